# Basic functions for import and adding of chrons
#
# On the "Tabelle1" sheet, the four rows that were tagged as "comp" in the
# type column (G) are reclassified as "add" rows, and the active
# selection is moved to G10 (reflecting the cell the user ended up on
# after making the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("G6").Value = "add"
$ws.Range("G7").Value = "add"
$ws.Range("G8").Value = "add"
$ws.Range("G9").Value = "add"

$ws.Range("G10").Select()
